# Auto-committed on 2022/03/28 週一
# Update the "CdSyndFee" DB layout sheet: rename the syndicated-loan-fee
# ("聯貸費用...") wording to corporate-finance-fee ("企金費用...") wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# D1 - table (Chinese) description
$ws.Range("D1").Value = "企金費用代碼檔"

# C9 - field description for SyndFeeCode
$ws.Range("C9").Value = "企金費用代碼"

# C10 - field description for SyndFeeItem
$ws.Range("C10").Value = "企金費用說明"

# Leave the cursor where the last edit was made.
$ws.Activate()
$ws.Range("C10").Select()
